$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "['Fc5', 'Cp5', 'F7', 'Ft7', 'T7', 'T9', 'T10', 'Tp7', 'O1', 'Iz']"
$ws.Range("C3").Value = "['Fc3', 'Fc1', 'Fcz', 'C5', 'C1', 'Cz', 'Fp1', 'Fpz', 'Fp2', 'Af7', 'Af3', 'Afz', 'F5', 'F3', 'F1', 'Fz', 'F2']"
$ws.Range("C4").Value = "['Fc2', 'Fc4', 'Fc6', 'C2', 'C4', 'C6', 'Af4', 'Af8', 'F4', 'F6', 'F8', 'Ft8', 'T8']"
$ws.Range("C5").Value = "['C3', 'Cp3', 'Cp1', 'P7', 'P5', 'P3', 'P1', 'Po7', 'Po3']"
$ws.Range("C6").Value = "['Cpz', 'Cp2', 'Cp4', 'Cp6', 'Tp8', 'Pz', 'P2', 'P4', 'P6', 'P8', 'Poz', 'Po4', 'Po8', 'Oz', 'O2']"

$ws.Range("B2").Value = 3
$ws.Range("B3").Value = 1
$ws.Range("B4").Value = 4
$ws.Range("B5").Value = 2
$ws.Range("B6").Value = 0
